$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume 1h (E) columns for rows 2-44, 47-51 ---
# Values that look like plain numbers ("1.005", "214.88", ...) are entered
# with a leading apostrophe to force Excel to keep them as text (matching the
# original inlineStr cell type), then the cell style is reset to "Normal" so
# no extra number-format / quote-prefix styling leaks into the cell.
$ws.Range("D2").Value = "25.925.10"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.641.01"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("D5").Value = "'214.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").Value = "'0.5048"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.91%  "
$ws.Range("D8").Value = "'0.2573"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "'0.06397"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").Value = "'19.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("D11").Value = "'0.07770"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "1.653.00"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").Value = "'4.259"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").Value = "1.868.89"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "'0.5433"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "0.0₅7914"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "'64.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("D18").Value = "25.959.10"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").Value = "'198.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").Value = "'4.376"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").Value = "'9.888"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").Value = "'5.968"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D24").Value = "'1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("D25").Value = "'1.880"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.47%  "
$ws.Range("D26").Value = "'141.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").Value = "'0.1134"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.42%  "
$ws.Range("D28").Value = "'6.818"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").Value = "'15.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "'1.239"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").Value = "'0.04936"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("D35").Value = "'2.365"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.79%  "
$ws.Range("D36").Value = "'0.8925"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("D37").Value = "'2.612"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("D38").Value = "1.143.18"
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("D39").Value = "'0.5549"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("D40").Value = "'0.01568"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").Value = "'1.003"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("D42").Value = "'5.720"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("D43").Value = "'0.8099"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").Value = "'99.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D47").Value = "'0.4522"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").Value = "'1.004"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").Value = "'54.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "'0.05063"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").Value = "'1.004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.45%  "

# --- Rows 32-34: only the Volume(1h) column changes, price stays the same ---
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("E34").Value = "  +0.06%  "

# --- Rows 45 & 46 swap coin identity (RocketPoolETH <-> BabyDogeCoin) ---
# along with updated price/volume figures
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₈120"
$ws.Range("E45").Value = "  +4.98%  "

$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.779.83"
$ws.Range("E46").Value = "  +0.62%  "
